# Applies the "6.0.0" release refresh to the FHIR StructureDefinition
# spreadsheet:
#   - Metadata sheet: bump Version, refresh publication Date, set the
#     Publisher display text, replace the stray duplicate "Contact" row
#     with a proper "Jurisdiction" row, and drop the leftover duplicate
#     row entirely.
#   - Elements sheet: give the root Extension row a Short/Definition that
#     describes this specific extension instead of the generic
#     "Extension" / "An Extension" boilerplate.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# The old row 11 duplicated row 10's "Contact" / "No display for
# ContactDetail" pair; remove it so everything below shifts up one row
# (A1:B21 -> A1:B20).
$meta.Rows("11").Delete()

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Claim Received Date"
$elements.Range("L2").Value = "Date the claim was received"
